$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.851.82"
$ws.Range("E2").Value = "  +5.34%  "
$ws.Range("D3").Value = "2.232.34"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'231.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").Value = "'0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "'61.71"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.18%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D10").Value = "'59.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  +5.48%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "2.562.63"
$ws.Range("D14").Value = "'15.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "'22.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "'5.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "2.236.21"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("D19").Value = "41.801.71"
$ws.Range("E19").Value = "  +5.32%  "
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").Value = "'72.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").Value = "'250.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "'9.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("D29").Value = "'167.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").Value = "'20.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("D32").Value = "'2.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'5.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.65%  "
$ws.Range("E35").Value = "  +3.21%  "
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("E37").Value = "  -4.65%  "
$ws.Range("D38").Value = "'3.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.70%  "
$ws.Range("D39").Value = "'2.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "'0.000258"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +32.14%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "'4.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("E43").Value = "  +4.45%  "
$ws.Range("D44").Value = "'8.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.48%  "
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").Value = "'0.0976"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.14%  "
$ws.Range("D47").Value = "'99.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("D48").Value = "1.480.63"
$ws.Range("D49").Value = "'16.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.05%  "
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "'52.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.18%  "
